$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: spaces -> underscores (columns A1:M1) ---
# (Order matches the original authoring order so new shared-string entries
#  land at the same indices as the reference edit.)
$ws.Range("A1").Value = "Last_Name "
$ws.Range("B1").Value = "First_Name"
$ws.Range("H1").Value = "Time_of_Experiment"
$ws.Range("G1").Value = "Date_of_Experiment "
$ws.Range("I1").Value = "Passage_Title "
$ws.Range("J1").Value = "Correct_Answer"
$ws.Range("K1").Value = "Incorrect_Answer"
$ws.Range("L1").Value = "No_Answer"
$ws.Range("M1").Value = "Total_Questions"

# --- New "Bad_Data " column header ---
$ws.Range("N1").Value = "Bad_Data "

# --- Update "Time of Experiment" (column H) values for every data row,
#     and populate the new "Bad_Data" boolean column (TRUE on rows 8 & 14) ---
# Each tuple is: row, new H value (time id), Bad_Data flag
$rowData = @(
    @(2,  160122, $false),
    @(3,  150255, $false),
    @(4,  143937, $false),
    @(5,  185654, $false),
    @(6,  191010, $false),
    @(7,  193919, $false),
    @(8,  194937, $true),
    @(9,  192957, $false),
    @(10, 194603, $false),
    @(11, 202514, $false),
    @(12, 204405, $false),
    @(13, 170520, $false),
    @(14, 192040, $true),
    @(15, 175032, $false),
    @(16, 175948, $false),
    @(17, 180053, $false),
    @(18, 181931, $false),
    @(19, 190525, $false)
)

foreach ($entry in $rowData) {
    $row = $entry[0]
    $timeValue = $entry[1]
    $badDataFlag = $entry[2]
    $ws.Range("H$row").Value = $timeValue
    $ws.Range("N$row").Value = $badDataFlag
}

# --- Refresh the view state to match the saved workbook (cosmetic) ---
$ws.Range("H15").Select() | Out-Null
